$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1) Remove the duplicated "RESERVA HIDRO" row (old row 10). This
#    shifts every row below it up by one and updates merged ranges
#    and the sheet dimension automatically.
# ---------------------------------------------------------------
$ws.Rows.Item(10).Delete()

# ---------------------------------------------------------------
# 2) Fix the order of the two header rows that used to be 8 & 9:
#    row 8 must now read "RESERVA PROGRAMADA A 50Hz PARA RPF"
#    row 9 must now read "RESERVA HIDRO"
# ---------------------------------------------------------------
$ws.Range("A8").Value = "RESERVA PROGRAMADA A 50Hz PARA RPF"
$ws.Range("A9").Value = "RESERVA HIDRO"

# ---------------------------------------------------------------
# 3) Row 7 ("RESERVA ROTANTE DEL PARQUE REGULANTE") becomes a
#    border-less, centered, merged A7:E7 cell (previously bordered
#    A7:F7).
# ---------------------------------------------------------------
$ws.Range("A7:F7").UnMerge()
$ws.Range("A7:F7").Clear()
$ws.Range("A7:E7").Merge()
$ws.Range("A7").Borders.LineStyle = -4142
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("A7").Value = "RESERVA ROTANTE DEL PARQUE REGULANTE"

# ---------------------------------------------------------------
# 4) Row 8 becomes a bordered, centered, merged A8:F8 cell - reuse
#    the existing bordered style by copying formats from A3:F3.
# ---------------------------------------------------------------
$ws.Range("A8:F8").Merge()
$ws.Range("A3:F3").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A8").Value = "RESERVA PROGRAMADA A 50Hz PARA RPF"

# ---------------------------------------------------------------
# 5) Rows 9-14 get split label/value merges.
# ---------------------------------------------------------------
$ws.Range("A9:C9").Merge()
$ws.Range("D9:F9").Merge()
$ws.Range("A10:C10").Merge()
$ws.Range("D10:F10").Merge()
$ws.Range("A11:C11").Merge()
$ws.Range("D11:F11").Merge()
$ws.Range("A12:E12").Merge()
$ws.Range("A13:E13").Merge()
$ws.Range("A14:E14").Merge()

# ---------------------------------------------------------------
# 6) Row 15 ("POTENCIA OPERABLE EN EL PARQUE REGULANTE") becomes a
#    border-less, centered, merged A15:F15 cell.
# ---------------------------------------------------------------
$ws.Range("A15:F15").UnMerge()
$ws.Range("A15:F15").Clear()
$ws.Range("A15:F15").Merge()
$ws.Range("A15").Borders.LineStyle = -4142
$ws.Range("A15").HorizontalAlignment = -4108
$ws.Range("A15").Value = "POTENCIA OPERABLE EN EL PARQUE REGULANTE"

# ---------------------------------------------------------------
# 7) Rows 16-19 get split label/value merges.
# ---------------------------------------------------------------
$ws.Range("A16:C16").Merge()
$ws.Range("D16:F16").Merge()
$ws.Range("A17:C17").Merge()
$ws.Range("D17:F17").Merge()
$ws.Range("A18:C18").Merge()
$ws.Range("D18:F18").Merge()
$ws.Range("A19:C19").Merge()
$ws.Range("D19:F19").Merge()

# ---------------------------------------------------------------
# 8) Rows 21-24 get split label/value merges (row 20 is untouched,
#    it already carries the bordered/merged/centered style).
# ---------------------------------------------------------------
$ws.Range("A21:C21").Merge()
$ws.Range("D21:F21").Merge()
$ws.Range("A22:C22").Merge()
$ws.Range("D22:F22").Merge()
$ws.Range("A23:C23").Merge()
$ws.Range("D23:F23").Merge()
$ws.Range("A24:C24").Merge()
$ws.Range("D24:F24").Merge()
